$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.222.33'
$ws.Cells.Item(2, 5).Value = '  +0.01%  '

$ws.Cells.Item(3, 4).Value = '1.889.91'
$ws.Cells.Item(3, 5).Value = '  -0.78%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(4, 5).Value = '  +0.32%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '322.66'
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 5).Value = '  -2.99%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 5).Value = '  +0.25%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4721'
$ws.Cells.Item(7, 4).NumberFormat = "General"
$ws.Cells.Item(7, 5).Value = '  +2.28%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4033'
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(8, 5).Value = '  -2.42%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '47.34'
$ws.Cells.Item(9, 4).NumberFormat = "General"
$ws.Cells.Item(9, 5).Value = '  -1.14%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.08010'
$ws.Cells.Item(10, 4).NumberFormat = "General"
$ws.Cells.Item(10, 5).Value = '  -0.35%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.9949'
$ws.Cells.Item(11, 4).NumberFormat = "General"
$ws.Cells.Item(11, 5).Value = '  -2.00%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '22.85'
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 5).Value = '  +2.96%  '

$ws.Cells.Item(13, 4).Value = '1.896.82'
$ws.Cells.Item(13, 5).Value = '  -0.03%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.914'
$ws.Cells.Item(14, 4).NumberFormat = "General"
$ws.Cells.Item(14, 5).Value = '  -0.66%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.036'
$ws.Cells.Item(15, 4).NumberFormat = "General"
$ws.Cells.Item(15, 5).Value = '  -1.32%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '89.23'
$ws.Cells.Item(16, 4).NumberFormat = "General"
$ws.Cells.Item(16, 5).Value = '  -0.10%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.003'
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 5).Value = '  +0.41%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.06631'
$ws.Cells.Item(18, 4).NumberFormat = "General"
$ws.Cells.Item(18, 5).Value = '  +0.78%  '

$ws.Cells.Item(19, 5).Value = '  -0.82%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.48'
$ws.Cells.Item(20, 4).NumberFormat = "General"
$ws.Cells.Item(20, 5).Value = '  -1.08%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.000'
$ws.Cells.Item(21, 4).NumberFormat = "General"
$ws.Cells.Item(21, 5).Value = '  -0.04%  '

$ws.Cells.Item(22, 4).Value = '29.248.58'
$ws.Cells.Item(22, 5).Value = '  +0.20%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.486'
$ws.Cells.Item(23, 4).NumberFormat = "General"
$ws.Cells.Item(23, 5).Value = '  -0.41%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.69'
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 5).Value = '  +2.40%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.176'
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 5).Value = '  -0.86%  '

$ws.Cells.Item(26, 4).Value = '2.191.93'
$ws.Cells.Item(26, 5).Value = '  +3.31%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '155.26'
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 5).Value = '  -1.29%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '19.63'
$ws.Cells.Item(28, 4).NumberFormat = "General"
$ws.Cells.Item(28, 5).Value = '  -0.79%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.922'
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 5).Value = '  +4.65%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.078'
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 5).Value = '  -2.43%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '117.10'
$ws.Cells.Item(31, 4).NumberFormat = "General"
$ws.Cells.Item(31, 5).Value = '  -0.19%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.021'
$ws.Cells.Item(32, 4).NumberFormat = "General"
$ws.Cells.Item(32, 5).Value = '  -2.34%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.09413'
$ws.Cells.Item(33, 4).NumberFormat = "General"
$ws.Cells.Item(33, 5).Value = '  -0.07%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.349'
$ws.Cells.Item(36, 4).NumberFormat = "General"
$ws.Cells.Item(36, 5).Value = '  -0.39%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02241'
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 5).Value = '  -0.46%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.06042'
$ws.Cells.Item(38, 4).NumberFormat = "General"
$ws.Cells.Item(38, 5).Value = '  -1.07%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.170'
$ws.Cells.Item(39, 4).NumberFormat = "General"
$ws.Cells.Item(39, 5).Value = '  -0.70%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '8.004'
$ws.Cells.Item(40, 4).NumberFormat = "General"
$ws.Cells.Item(40, 5).Value = '  -5.47%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.5816'
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 5).Value = '  -0.86%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.1826'
$ws.Cells.Item(42, 4).NumberFormat = "General"
$ws.Cells.Item(42, 5).Value = '  -0.11%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '10.01'
$ws.Cells.Item(43, 4).NumberFormat = "General"
$ws.Cells.Item(43, 5).Value = '  -1.53%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.284'
$ws.Cells.Item(44, 4).NumberFormat = "General"
$ws.Cells.Item(44, 5).Value = '  +3.66%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.07702'
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 5).Value = '  +2.62%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.366'
$ws.Cells.Item(46, 4).NumberFormat = "General"
$ws.Cells.Item(46, 5).Value = '  +0.55%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '12.16'
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 5).Value = '  -0.24%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.5470'
$ws.Cells.Item(48, 4).NumberFormat = "General"
$ws.Cells.Item(48, 5).Value = '  -1.56%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.905'
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 5).Value = '  -1.06%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '113.17'
$ws.Cells.Item(50, 4).NumberFormat = "General"
$ws.Cells.Item(50, 5).Value = '  -0.12%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.2957'
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 5).Value = '  +1.04%  '

# Row 34/35 swap: ARBITRUM <-> HuobiToken (full row content swap incl. B, C, D, E)
$ws.Cells.Item(34, 2).Value = 'HuobiToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.530'
$ws.Cells.Item(34, 4).NumberFormat = "General"
$ws.Cells.Item(34, 5).Value = '  -0.43%  '

$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.378'
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 5).Value = '  -3.55%  '
